$wb = $excel.ActiveWorkbook

# --- Add two new driver worksheets at the front of the workbook ---
# Final desired order: colrow2, headerrow2, blank, standard, headeronly
# Insert "headerrow2" before the current first sheet, then insert "colrow2"
# before "headerrow2" so that colrow2 ends up first, headerrow2 second.
$firstSheet = $wb.Worksheets.Item(1)
$newHeaderRow2 = $wb.Worksheets.Add($firstSheet)
$newHeaderRow2.Name = "headerrow2"
$newColRow2 = $wb.Worksheets.Add($newHeaderRow2)
$newColRow2.Name = "colrow2"

# Re-fetch the worksheet objects by name so that Range()/Activate() calls
# below are bound to the correct sheet.
$colrow2 = $wb.Worksheets.Item("colrow2")
$headerrow2 = $wb.Worksheets.Item("headerrow2")

# --- "colrow2": same table as "standard" but data starts one column over ---
# (columns B:C instead of A:B)
$colrow2.Range("B1").Value = "uid"
$colrow2.Range("C1").Value = "uname"
$colrow2.Range("B2").Value = 1
$colrow2.Range("C2").Value = "mike"
$colrow2.Range("B3").Value = 2
$colrow2.Range("B4").Value = 3
$colrow2.Range("C4").Value = "henry"

# --- "headerrow2": same table as "standard" but data starts one row down ---
# (rows 2:5 instead of 1:4)
$headerrow2.Range("A2").Value = "uid"
$headerrow2.Range("B2").Value = "uname"
$headerrow2.Range("A3").Value = 1
$headerrow2.Range("B3").Value = "mike"
$headerrow2.Range("A4").Value = 2
$headerrow2.Range("A5").Value = 3
$headerrow2.Range("B5").Value = "henry"

# Select the header row (A1:XFD1) on "headerrow2", mirroring the target file.
$headerrow2.Range("A1:XFD1").Select() | Out-Null

# Make "colrow2" the active (first / tabSelected) sheet again.
$colrow2.Activate() | Out-Null
